# Romania_FX.xlsx update:
#  - row 220 (2023-08-01): correct the high/close values
#  - append three new monthly rows (221-223) with the same layout/style
#    as the existing data rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct existing row 220 ---
$ws.Cells.Item(220, 4).Value = 4.58647   # D220 high
$ws.Cells.Item(220, 6).Value = 4.5542    # F220 close

# --- New rows to append ---
$newRows = @(
    @{ Row = 221; A = 45170.33333333334; C = 4.5546;  D = 4.7417;  E = 4.54653; F = 4.7021; G = 0 },
    @{ Row = 222; A = 45201.375;         C = 4.6964;  D = 4.7602;  E = 4.64967; F = 4.6949; G = 0 },
    @{ Row = 223; A = 45231.375;         C = 4.6941;  D = 4.7224;  E = 4.61965; F = 4.6438; G = 0 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Match the formatting of the last existing data row, column by column
    for ($col = 1; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).NumberFormat = $ws.Cells.Item(220, $col).NumberFormat
    }
    $ws.Cells.Item($row, 1).Font.Bold = $true
    $ws.Cells.Item($row, 1).HorizontalAlignment = -4108   # xlCenter
    $ws.Cells.Item($row, 1).VerticalAlignment = -4160     # xlTop

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = "FX_IDC:USDRON"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}
